$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value2 = 32.533333
$ws.Range("I11").Value2 = 32.533333
$ws.Range("K11").Value2 = 32.533333
$ws.Range("M11").Value2 = 107.466667
$ws.Range("H113").Value2 = 2999.5
$ws.Range("I113").Value2 = 2999.5
$ws.Range("K113").Value2 = 2999.5
$ws.Range("M113").Value2 = 254.5
$ws.Range("H132").Value2 = 53772.58
$ws.Range("I132").Value2 = 54653.637
$ws.Range("K132").Value2 = 163960.911
$ws.Range("M132").Value2 = -161430.911
$ws.Range("H137").Value2 = 1311636.4
$ws.Range("I137").Value2 = 869373.6
$ws.Range("K137").Value2 = 2608120.8
$ws.Range("M137").Value2 = -2605570.8
$ws.Range("H138").Value2 = 2576.8845
$ws.Range("I138").Value2 = 2050.4
$ws.Range("K138").Value2 = 6151.200000000001
$ws.Range("M138").Value2 = -1011.200000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 5956180
$ws.Range("J32").Value2 = 40659
$ws.Range("L32").Value2 = 40659
$ws.Range("N32").Value2 = -41233
$ws.Range("H45").Value2 = 3634.3215
$ws.Range("I45").Value2 = 3773.0588
$ws.Range("J45").Value2 = 3419.9092
$ws.Range("K45").Value2 = 3773.0588
$ws.Range("L45").Value2 = 3419.9092
$ws.Range("M45").Value2 = -3396.0588
$ws.Range("N45").Value2 = -4173.9092
$ws.Range("H61").Value2 = 1973226.8
$ws.Range("I61").Value2 = 2579046
$ws.Range("K61").Value2 = 2579046
$ws.Range("M61").Value2 = -2578834
$ws.Range("H63").Value2 = 7023.3076
$ws.Range("I63").Value2 = 2099.8
$ws.Range("J63").Value2 = 10100.5
$ws.Range("K63").Value2 = 2099.8
$ws.Range("L63").Value2 = 10100.5
$ws.Range("M63").Value2 = -1413.8
$ws.Range("N63").Value2 = -11472.5
$ws.Range("H66").Value2 = 7023.3076
$ws.Range("I66").Value2 = 2099.8
$ws.Range("J66").Value2 = 10100.5
$ws.Range("K66").Value2 = 10499
$ws.Range("L66").Value2 = 50502.5
$ws.Range("M66").Value2 = -7067
$ws.Range("N66").Value2 = -57366.5
$ws.Range("H74").Value2 = 2909936
$ws.Range("I74").Value2 = 3380687.2
$ws.Range("J74").Value2 = 6970
$ws.Range("K74").Value2 = 3380687.2
$ws.Range("L74").Value2 = 6970
$ws.Range("M74").Value2 = -3379813.2
$ws.Range("N74").Value2 = -8718
$ws.Range("H77").Value2 = 2909936
$ws.Range("I77").Value2 = 3380687.2
$ws.Range("J77").Value2 = 6970
$ws.Range("K77").Value2 = 16903436
$ws.Range("L77").Value2 = 34850
$ws.Range("M77").Value2 = -16899068
$ws.Range("N77").Value2 = -43586
$ws.Range("H80").Value2 = 77115.75
$ws.Range("J80").Value2 = 77821
$ws.Range("L80").Value2 = 77821
$ws.Range("N80").Value2 = -79817
$ws.Range("H83").Value2 = 77115.75
$ws.Range("J83").Value2 = 77821
$ws.Range("L83").Value2 = 233463
$ws.Range("N83").Value2 = -243447
$ws.Range("H132").Value2 = 708073.6
$ws.Range("I132").Value2 = 777384.3
$ws.Range("K132").Value2 = 2332152.9
$ws.Range("M132").Value2 = -2329622.9
$ws.Range("H136").Value2 = 1973226.8
$ws.Range("I136").Value2 = 2579046
$ws.Range("K136").Value2 = 7737138
$ws.Range("M136").Value2 = -7734588

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 1656.7037
$ws.Range("I20").Value2 = 1781.5
$ws.Range("J20").Value2 = 1407.1111
$ws.Range("K20").Value2 = 1781.5
$ws.Range("L20").Value2 = 1407.1111
$ws.Range("M20").Value2 = -1534.5
$ws.Range("N20").Value2 = -1901.1111
$ws.Range("H82").Value2 = 52820.734
$ws.Range("I82").Value2 = 23100
$ws.Range("J82").Value2 = 63628.273
$ws.Range("K82").Value2 = 23100
$ws.Range("L82").Value2 = 63628.273
$ws.Range("M82").Value2 = -22717
$ws.Range("N82").Value2 = -64394.273
$ws.Range("H85").Value2 = 52820.734
$ws.Range("I85").Value2 = 23100
$ws.Range("J85").Value2 = 63628.273
$ws.Range("K85").Value2 = 23100
$ws.Range("L85").Value2 = 63628.273
$ws.Range("M85").Value2 = -21774
$ws.Range("N85").Value2 = -66280.273
$ws.Range("H99").Value2 = 14760.823
$ws.Range("I99").Value2 = 16425.953
$ws.Range("K99").Value2 = 16425.953
$ws.Range("M99").Value2 = -14927.953
$ws.Range("H134").Value2 = 1244951.2
$ws.Range("I134").Value2 = 1461478.2
$ws.Range("K134").Value2 = 4384434.6
$ws.Range("M134").Value2 = -4381899.6

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 914.9231
$ws.Range("I16").Value2 = 788.2
$ws.Range("J16").Value2 = 1337.3334
$ws.Range("K16").Value2 = 788.2
$ws.Range("L16").Value2 = 1337.3334
$ws.Range("M16").Value2 = -501.2
$ws.Range("N16").Value2 = -1911.3334
$ws.Range("H31").Value2 = 4219.1313
$ws.Range("I31").Value2 = 1503.4706
$ws.Range("K31").Value2 = 1503.4706
$ws.Range("M31").Value2 = -1208.4706
$ws.Range("H34").Value2 = 4219.1313
$ws.Range("I34").Value2 = 1503.4706
$ws.Range("K34").Value2 = 1503.4706
$ws.Range("M34").Value2 = -1301.4706
$ws.Range("H52").Value2 = 0
$ws.Range("I52").Value2 = 0
$ws.Range("K52").Value2 = 0
$ws.Range("M52").ClearContents()
$ws.Range("H58").Value2 = 1375327.6
$ws.Range("I58").Value2 = 2059158.4
$ws.Range("K58").Value2 = 2059158.4
$ws.Range("M58").Value2 = -2058955.4
$ws.Range("H68").Value2 = 80255.69500000001
$ws.Range("I68").Value2 = 70890
$ws.Range("J68").Value2 = 81036.164
$ws.Range("K68").Value2 = 70890
$ws.Range("L68").Value2 = 81036.164
$ws.Range("M68").Value2 = -70141
$ws.Range("N68").Value2 = -82534.164
$ws.Range("H71").Value2 = 80255.69500000001
$ws.Range("I71").Value2 = 70890
$ws.Range("J71").Value2 = 81036.164
$ws.Range("K71").Value2 = 212670
$ws.Range("L71").Value2 = 243108.492
$ws.Range("M71").Value2 = -208926
$ws.Range("N71").Value2 = -250596.492
$ws.Range("H113").Value2 = 914.9231
$ws.Range("I113").Value2 = 788.2
$ws.Range("J113").Value2 = 1337.3334
$ws.Range("K113").Value2 = 788.2
$ws.Range("L113").Value2 = 1337.3334
$ws.Range("M113").Value2 = 1381.8
$ws.Range("N113").Value2 = -5677.3334
$ws.Range("H122").Value2 = 2544.6
$ws.Range("I122").Value2 = 1032.625
$ws.Range("K122").Value2 = 3097.875
$ws.Range("M122").Value2 = -647.875
$ws.Range("H125").Value2 = 40000
$ws.Range("J125").Value2 = 40000
$ws.Range("L125").Value2 = 40000
$ws.Range("N125").Value2 = -44920
$ws.Range("H132").Value2 = 9274507
$ws.Range("I132").Value2 = 16912.334
$ws.Range("K132").Value2 = 50737.00199999999
$ws.Range("M132").Value2 = -48207.00199999999
$ws.Range("H136").Value2 = 1375327.6
$ws.Range("I136").Value2 = 2059158.4
$ws.Range("K136").Value2 = 6177475.199999999
$ws.Range("M136").Value2 = -6174925.199999999
$ws.Range("H137").Value2 = 70496
$ws.Range("J137").Value2 = 69999
$ws.Range("L137").Value2 = 69999
$ws.Range("N137").Value2 = -80199

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value2 = 498.8889
$ws.Range("I34").Value2 = 500
$ws.Range("J34").Value2 = 498.57144
$ws.Range("K34").Value2 = 1500
$ws.Range("L34").Value2 = 1495.71432
$ws.Range("M34").Value2 = -1416
$ws.Range("N34").Value2 = -1663.71432
$ws.Range("H39").Value2 = 2400
$ws.Range("J39").Value2 = 2000
$ws.Range("L39").Value2 = 6000
$ws.Range("N39").Value2 = -6588
$ws.Range("H55").Value2 = 1741
$ws.Range("J55").Value2 = 995
$ws.Range("L55").Value2 = 2985
$ws.Range("N55").Value2 = -3339
$ws.Range("H113").Value2 = 2008.7142
$ws.Range("J113").Value2 = 1937.3
$ws.Range("L113").Value2 = 5811.9
$ws.Range("N113").Value2 = -10151.9
$ws.Range("H114").Value2 = 1367.7407
$ws.Range("J114").Value2 = 3482.1
$ws.Range("L114").Value2 = 10446.3
$ws.Range("N114").Value2 = -16954.3

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value2 = 76650
$ws.Range("J123").Value2 = 76650
$ws.Range("L123").Value2 = 76650
$ws.Range("N123").Value2 = -81550

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value2 = 2480.3333
$ws.Range("I46").Value2 = 1249.25
$ws.Range("J46").Value2 = 2832.0715
$ws.Range("K46").Value2 = 1249.25
$ws.Range("L46").Value2 = 2832.0715
$ws.Range("M46").Value2 = -1061.25
$ws.Range("N46").Value2 = -3208.0715
$ws.Range("H55").Value2 = 1459.25
$ws.Range("I55").Value2 = 175.5
$ws.Range("J55").Value2 = 1887.1666
$ws.Range("K55").Value2 = 175.5
$ws.Range("L55").Value2 = 1887.1666
$ws.Range("M55").Value2 = -2.5
$ws.Range("N55").Value2 = -2233.1666
$ws.Range("H100").Value2 = 9791
$ws.Range("I100").Value2 = 2806.6
$ws.Range("J100").Value2 = 27252
$ws.Range("K100").Value2 = 2806.6
$ws.Range("L100").Value2 = 27252
$ws.Range("M100").Value2 = -2265.6
$ws.Range("N100").Value2 = -28334
$ws.Range("H132").Value2 = 846431.2
$ws.Range("I132").Value2 = 1019475.9
$ws.Range("J132").Value2 = 5928.4287
$ws.Range("K132").Value2 = 3058427.7
$ws.Range("L132").Value2 = 17785.2861
$ws.Range("M132").Value2 = -3055897.7
$ws.Range("N132").Value2 = -22845.2861

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H129").Value2 = 53999
$ws.Range("J129").Value2 = 53999
$ws.Range("L129").Value2 = 53999
$ws.Range("N129").Value2 = -63999
$ws.Range("H136").Value2 = 8471712
$ws.Range("I136").Value2 = 9774278
$ws.Range("J136").Value2 = 5033.1665
$ws.Range("K136").Value2 = 29322834
$ws.Range("L136").Value2 = 15099.4995
$ws.Range("M136").Value2 = -29320284
$ws.Range("N136").Value2 = -20199.4995
